$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{Cell = 'D2'; Value = '302.79'},
    @{Cell = 'E2'; Value = '-4.49%'},
    @{Cell = 'D3'; Value = '35.30'},
    @{Cell = 'E3'; Value = '-2.27%'},
    @{Cell = 'D4'; Value = '5.070'},
    @{Cell = 'E4'; Value = '-1.88%'},
    @{Cell = 'D5'; Value = '0.08004'},
    @{Cell = 'E5'; Value = '-2.96%'},
    @{Cell = 'D6'; Value = '1.938'},
    @{Cell = 'E6'; Value = '-9.79%'},
    @{Cell = 'D7'; Value = '4.064'},
    @{Cell = 'E7'; Value = '-1.92%'},
    @{Cell = 'D8'; Value = '7.766'},
    @{Cell = 'E8'; Value = '-3.43%'},
    @{Cell = 'D9'; Value = '2.927'},
    @{Cell = 'E9'; Value = '4.50%'},
    @{Cell = 'D10'; Value = '0.9247'},
    @{Cell = 'E10'; Value = '-0.31%'},
    @{Cell = 'D11'; Value = '0.1230'},
    @{Cell = 'E11'; Value = '20.54%'},
    @{Cell = 'D12'; Value = '0.1859'},
    @{Cell = 'E12'; Value = '-1.51%'},
    @{Cell = 'D13'; Value = '0.09643'},
    @{Cell = 'E13'; Value = '2.94%'},
    @{Cell = 'D14'; Value = '0.03626'},
    @{Cell = 'E14'; Value = '1.08%'},
    @{Cell = 'D15'; Value = '0.09853'},
    @{Cell = 'E15'; Value = '-0.73%'},
    @{Cell = 'D16'; Value = '0.001386'},
    @{Cell = 'E16'; Value = '-3.44%'},
    @{Cell = 'D17'; Value = '0.005771'},
    @{Cell = 'E17'; Value = '2.08%'},
    @{Cell = 'D18'; Value = '3.506'},
    @{Cell = 'E18'; Value = '1.13%'},
    @{Cell = 'E19'; Value = '1.04%'},
    @{Cell = 'D20'; Value = '0.1310'},
    @{Cell = 'E20'; Value = '-1.58%'},
    @{Cell = 'D21'; Value = '5.056'},
    @{Cell = 'E21'; Value = '-2.41%'},
    @{Cell = 'E22'; Value = '12.43%'},
    @{Cell = 'D23'; Value = '0.04530'},
    @{Cell = 'E23'; Value = '-1.57%'},
    @{Cell = 'E24'; Value = '-2.40%'},
    @{Cell = 'D25'; Value = '0.004831'},
    @{Cell = 'E25'; Value = '2.07%'},
    @{Cell = 'E26'; Value = '0.04%'},
    @{Cell = 'D27'; Value = '0.0003006'},
    @{Cell = 'E27'; Value = '-33.30%'},
    @{Cell = 'D39'; Value = '0.01921'},
    @{Cell = 'E39'; Value = '-4.04%'},
    @{Cell = 'D40'; Value = '0.04720'},
    @{Cell = 'E40'; Value = '-4.52%'},
    @{Cell = 'D41'; Value = '0.007549'},
    @{Cell = 'E41'; Value = '-4.88%'},
    @{Cell = 'D42'; Value = '0.009630'},
    @{Cell = 'E42'; Value = '22.74%'},
    @{Cell = 'D43'; Value = '0.1327'},
    @{Cell = 'E43'; Value = '-5.24%'},
    @{Cell = 'D44'; Value = '0.002113'},
    @{Cell = 'E44'; Value = '-0.01%'},
    @{Cell = 'D45'; Value = '0.01014'},
    @{Cell = 'E45'; Value = '-13.64%'},
    @{Cell = 'D46'; Value = '0.00006248'},
    @{Cell = 'E46'; Value = '-4.02%'},
    @{Cell = 'D47'; Value = '0.00000000751'},
    @{Cell = 'E48'; Value = '88.56%'},
    @{Cell = 'E49'; Value = '-21.72%'},
    @{Cell = 'D50'; Value = '0.00002102'},
    @{Cell = 'D51'; Value = '0.0002002'},
)

foreach ($update in $updates) {
    $c = $ws.Range($update.Cell)
    $origStyle = $c.Style
    $c.NumberFormat = "@"
    $c.Value = $update.Value
    $c.Style = $origStyle
}
